$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10679.471
$ws.Range("I32").Value = 8124.75
$ws.Range("J32").Value = 11465.538
$ws.Range("K32").Value = 8124.75
$ws.Range("L32").Value = 11465.538
$ws.Range("M32").Value = -7798.75
$ws.Range("N32").Value = -12117.538
$ws.Range("H64").Value = 18531826
$ws.Range("I64").Value = 33342872
$ws.Range("K64").Value = 33342872
$ws.Range("M64").Value = -33342624
$ws.Range("H67").Value = 18531826
$ws.Range("I67").Value = 33342872
$ws.Range("K67").Value = 33342872
$ws.Range("M67").Value = -33342014
$ws.Range("H74").Value = 14296688
$ws.Range("I74").Value = 20415982
$ws.Range("K74").Value = 20415982
$ws.Range("M74").Value = -20415046
$ws.Range("H77").Value = 14296688
$ws.Range("I77").Value = 20415982
$ws.Range("K77").Value = 102079910
$ws.Range("M77").Value = -102075230
$ws.Range("H80").Value = 373
$ws.Range("J80").Value = 298.875
$ws.Range("L80").Value = 896.625
$ws.Range("N80").Value = -2892.625
$ws.Range("H83").Value = 373
$ws.Range("J83").Value = 298.875
$ws.Range("L83").Value = 2689.875
$ws.Range("N83").Value = -12673.875
$ws.Range("H96").Value = 1038.25
$ws.Range("I96").Value = 718
$ws.Range("K96").Value = 2154
$ws.Range("M96").Value = -781
$ws.Range("H98").Value = 3230.9412
$ws.Range("I98").Value = 3043.6538
$ws.Range("K98").Value = 3043.6538
$ws.Range("M98").Value = -1545.6538
$ws.Range("H103").Value = 525
$ws.Range("I103").Value = 510.66666
$ws.Range("J103").Value = 532.1667
$ws.Range("K103").Value = 1531.99998
$ws.Range("L103").Value = 1596.5001
$ws.Range("M103").Value = -945.9999800000001
$ws.Range("N103").Value = -2768.5001
$ws.Range("H122").Value = 3230.9412
$ws.Range("I122").Value = 3043.6538
$ws.Range("K122").Value = 9130.9614
$ws.Range("M122").Value = -6680.9614
$ws.Range("H129").Value = 1860.8572
$ws.Range("I129").Value = 1500
$ws.Range("J129").Value = 2131.5
$ws.Range("K129").Value = 4500
$ws.Range("L129").Value = 6394.5
$ws.Range("M129").Value = 500
$ws.Range("N129").Value = -16394.5
$ws.Range("H131").Value = 1072.1
$ws.Range("I131").Value = 577.125
$ws.Range("K131").Value = 1731.375
$ws.Range("M131").Value = 3308.625
$ws.Range("H132").Value = 335417.44
$ws.Range("I132").Value = 393705.47
$ws.Range("K132").Value = 1181116.41
$ws.Range("M132").Value = -1178586.41
$ws.Range("H137").Value = 6639.357
$ws.Range("J137").Value = 8888.909
$ws.Range("L137").Value = 26666.727
$ws.Range("N137").Value = -31766.727
$ws.Range("H138").Value = 3115.7036
$ws.Range("J138").Value = 3955.1765
$ws.Range("L138").Value = 11865.5295
$ws.Range("N138").Value = -22145.5295

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2779629
$ws.Range("I32").Value = 1447.697
$ws.Range("J32").Value = 15876769
$ws.Range("K32").Value = 1447.697
$ws.Range("L32").Value = 15876769
$ws.Range("M32").Value = -1160.697
$ws.Range("N32").Value = -15877343
$ws.Range("H45").Value = 2305.5334
$ws.Range("J45").Value = 2873.25
$ws.Range("L45").Value = 2873.25
$ws.Range("N45").Value = -3627.25

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 890361.2
$ws.Range("I134").Value = 1113158.2
$ws.Range("J134").Value = 19427.182
$ws.Range("K134").Value = 3339474.6
$ws.Range("L134").Value = 58281.546
$ws.Range("M134").Value = -3336939.6
$ws.Range("N134").Value = -63351.546

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 18521724
$ws.Range("I16").Value = 33335022
$ws.Range("K16").Value = 33335022
$ws.Range("M16").Value = -33334735
$ws.Range("H94").Value = 4107.5
$ws.Range("I94").Value = 1359.2
$ws.Range("K94").Value = 1359.2
$ws.Range("M94").Value = -908.2
$ws.Range("H107").Value = 574.8333
$ws.Range("I107").Value = 574.8333
$ws.Range("K107").Value = 574.8333
$ws.Range("M107").Value = 1345.1667
$ws.Range("H113").Value = 18521724
$ws.Range("I113").Value = 33335022
$ws.Range("K113").Value = 33335022
$ws.Range("M113").Value = -33332852
$ws.Range("H132").Value = 34743.418
$ws.Range("I132").Value = 37220.09
$ws.Range("J132").Value = 7500
$ws.Range("K132").Value = 111660.27
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -109130.27
$ws.Range("N132").Value = -27560
$ws.Range("H139").Value = 95000
$ws.Range("I139").Value = 60000
$ws.Range("J139").Value = 130000
$ws.Range("K139").Value = 60000
$ws.Range("L139").Value = 130000
$ws.Range("M139").Value = -54860
$ws.Range("N139").Value = -140280

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 371.625
$ws.Range("J7").Value = 300
$ws.Range("L7").Value = 900
$ws.Range("N7").Value = -1124
$ws.Range("H14").Value = 569.125
$ws.Range("I14").Value = 569.125
$ws.Range("K14").Value = 1707.375
$ws.Range("M14").Value = -1534.375
$ws.Range("H50").Value = 1108.7858
$ws.Range("I50").Value = 1144.4166
$ws.Range("K50").Value = 3433.2498
$ws.Range("M50").Value = -2952.2498
$ws.Range("H53").Value = 1108.7858
$ws.Range("I53").Value = 1144.4166
$ws.Range("K53").Value = 3433.2498
$ws.Range("M53").Value = -2952.2498
$ws.Range("H92").Value = 190
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H107").Value = 3898.1487
$ws.Range("I107").Value = 584.3333
$ws.Range("J107").Value = 4740.644
$ws.Range("K107").Value = 1752.9999
$ws.Range("L107").Value = 14221.932
$ws.Range("M107").Value = 167.0001
$ws.Range("N107").Value = -18061.932
$ws.Range("H117").Value = 4233.1055
$ws.Range("J117").Value = 3503.6428
$ws.Range("L117").Value = 10510.9284
$ws.Range("N117").Value = -17394.9284
$ws.Range("H133").Value = 5991.25
$ws.Range("I133").Value = 6204.2856
$ws.Range("J133").Value = 4500
$ws.Range("K133").Value = 18612.8568
$ws.Range("L133").Value = 13500
$ws.Range("M133").Value = -13552.8568
$ws.Range("N133").Value = -23620
$ws.Range("H137").Value = 11722.143
$ws.Range("I137").Value = 10009.167
$ws.Range("J137").Value = 22000
$ws.Range("K137").Value = 30027.501
$ws.Range("L137").Value = 66000
$ws.Range("M137").Value = -24927.501
$ws.Range("N137").Value = -76200

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 700
$ws.Range("J21").Value = 700
$ws.Range("L21").Value = 700
$ws.Range("N21").Value = -1170
$ws.Range("H102").Value = 4185.569
$ws.Range("I102").Value = 2195.5454
$ws.Range("J102").Value = 5695.241
$ws.Range("K102").Value = 2195.5454
$ws.Range("L102").Value = 5695.241
$ws.Range("M102").Value = -573.5454
$ws.Range("N102").Value = -8939.241
$ws.Range("H126").Value = 11911167
$ws.Range("I126").Value = 35716650
$ws.Range("J126").Value = 8427.357
$ws.Range("K126").Value = 107149950
$ws.Range("L126").Value = 25282.071
$ws.Range("M126").Value = -107147480
$ws.Range("N126").Value = -30222.071

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 2432.5
$ws.Range("I55").Value = 577.5333000000001
$ws.Range("K55").Value = 577.5333000000001
$ws.Range("M55").Value = -404.5333000000001
$ws.Range("H132").Value = 6415.077
$ws.Range("I132").Value = 7269.2
$ws.Range("K132").Value = 21807.6
$ws.Range("M132").Value = -19277.6

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 20761.875
$ws.Range("J26").Value = 20872.572
$ws.Range("L26").Value = 20872.572
$ws.Range("N26").Value = -21458.572
$ws.Range("H96").Value = 3120.842
$ws.Range("I96").Value = 2116.3333
$ws.Range("K96").Value = 2116.3333
$ws.Range("M96").Value = -743.3332999999998
